$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.028.68"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.674.75"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.20"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.46"
$ws.Range("E9").Value = "  +5.50%  "
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.911.58"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.676.55"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.29"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.023.86"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.19"
$ws.Range("E18").Value = "  +2.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "235.90"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.27"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.08"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.45"
$ws.Range("E27").Value = "  +3.56%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.545.87"
$ws.Range("E33").Value = "  +7.01%  "
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("E35").Value = "  +4.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("E40").Value = "  +4.42%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.76"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.55"
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.818.64"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.780"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.65"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.03"
$ws.Range("E50").Value = "  +6.25%  "
$ws.Range("E51").Value = "  +0.29%  "
